# contingencies with rene fine
#
# The line set grows from 6 to 8 lines (line7, line8 appended after line6),
# which pushes the existing "extr" contingency rows down by two rows.
# Two brand-new contingency rows (extr7, extr8) are appended at the end
# with freshly generated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: becomes line7 (was extr1) ---
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

# --- Row 9: becomes line8 (was extr2) ---
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $false

# --- Row 10: becomes extr1 (was extr3) ---
$ws.Cells.Item(10, 2).Value = "extr1"
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

# --- Row 11: becomes extr2 (was extr4) ---
$ws.Cells.Item(11, 2).Value = "extr2"
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = $true

# --- Row 12: becomes extr3 (was extr5) ---
$ws.Cells.Item(12, 2).Value = "extr3"
$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 4).Value = 11
$ws.Cells.Item(12, 5).Value = $true

# --- Row 13: becomes extr4 (was extr6) ---
$ws.Cells.Item(13, 2).Value = "extr4"
$ws.Cells.Item(13, 3).Value = 7
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = $false

# --- Row 14: becomes extr5 (was extr7) ---
$ws.Cells.Item(14, 2).Value = "extr5"
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = $false

# --- Row 15: becomes extr6 (was extr8) ---
$ws.Cells.Item(15, 2).Value = "extr6"
$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $false

# --- New row 16: extr7 ---
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "extr7"
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $true

# --- New row 17: extr8 ---
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "extr8"
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $true

# Column A in the data rows carries the bold/centered/bordered style used
# throughout the table (same as A2:A15) - copy it onto the two new rows.
$ws.Cells.Item(15, 1).Copy()
$ws.Cells.Item(16, 1).PasteSpecial(-4122)
$ws.Cells.Item(17, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
